$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New entry (Sr. No 4) mirrors the layout of the previous entry (row 32):
# copy its formatting down to the new row first, then fill in the data.
$ws.Range("A32:F32").Copy()
$ws.Range("A34:F34").PasteSpecial(-4122)

$ws.Range("A34").Value = 4
$ws.Range("B34").Value = 45187
$ws.Range("C34").Value = "SLH/1993"
$ws.Range("D34").Value = "Shree Laxmi Lighting Hub"
$ws.Range("E34").Value = 1432
$ws.Range("F34").Formula = "=E34"

# Move the active selection down past the newly added row, matching
# where the cursor would land after entering this new block of data.
$ws.Range("A35").Select() | Out-Null
